# Add a new "Team adviser" member (IBM Japan / Mamoru Kitagawa) to both the
# Japanese entry sheet and the English entry sheet, and leave the selection
# on H20 (matching the saved UI state in the target file).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("記入用紙")

# --- Japanese member list: new "チームアドバイザー" row (row 11) ---
# Copy the formatting of an existing data row (A7:D7) down onto the blank
# row so the new entry matches the look of the other member rows.
$ws.Range("A7:D7").Copy($ws.Range("A11:D11"))
$ws.Range("A11").Value = "日本アイ・ビー・エム株式会社"
$ws.Range("B11").Value = "テクニカル・エキスパート本部"
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = "北川　衛"

# Row 12 (still blank) also picks up the same row formatting.
$ws.Range("A7:D7").Copy($ws.Range("A12:D12"))
$ws.Range("A12").ClearContents()
$ws.Range("B12").ClearContents()
$ws.Range("C12").ClearContents()
$ws.Range("D12").ClearContents()

# --- English member list: new "Team adviser" row (row 23) ---
$ws.Range("A7:D7").Copy($ws.Range("A23:D23"))
$ws.Range("A23").Value = "IBM Japan Co., Ltd"
$ws.Range("B23").ClearContents()
$ws.Range("C23").Value = ""
$ws.Range("D23").Value = "Mamoru Kitaagwa"

# The "Leader" row company name (A19) was re-entered by the author; keep the
# same visible text.
$ws.Range("A19").Value = "Sumitomo Mitsui Trust Systems & Services Co., Ltd"

# Restore the saved selection (H20) as recorded in the workbook.
$ws.Range("H20").Select()

Write-Output "done"
